$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 3 (columns 6-8): "41 590 000,00" -> "41 900 000,00"
for ($c = 6; $c -le 8; $c++) {
    $cell = $tbl.Cell(3, $c)
    if ($cell.Range.Text.Contains("41 590 000,00")) {
        $cell.Range.Text = "41 900 000,00"
    }
}

# Row 4 (columns 6-8): "3 743 100,00" -> "3 771 000,00"
for ($c = 6; $c -le 8; $c++) {
    $cell = $tbl.Cell(4, $c)
    if ($cell.Range.Text.Contains("3 743 100,00")) {
        $cell.Range.Text = "3 771 000,00"
    }
}
